# Add 2022-Q1 fund-holding data: rename the existing "总计" (Total) sheet
# to "2022-Q1" and populate it with the new quarter's fund detail table,
# then append a fresh "总计" sheet (cloned from the pre-edit Total sheet)
# that rolls up the new quarter alongside the historical ones.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: clone the current "总计" sheet (sheet index 4) BEFORE it gets
# overwritten. The clone becomes the new roll-up "总计" sheet; it starts
# out as an exact copy (same formatting/structure) of the original.
# ---------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item(4)
$oldTotal.Copy($null, $oldTotal)
$newTotal = $wb.Worksheets.Item(5)

# Rename the original sheet to "2022-Q1" (it will hold the new quarter's
# fund-level detail) and rename the clone back to "总计".
$oldTotal.Name = "2022-Q1"
$newTotal.Name = "总计"

$detail = $wb.Worksheets.Item(4)
$total = $wb.Worksheets.Item(5)

# ---------------------------------------------------------------------
# Step 2: rebuild the "2022-Q1" sheet with the fund-level holding detail.
# ---------------------------------------------------------------------
$detail.UsedRange.Clear()

$fundRows = @(
    @("110025", "易方达资源行业混合", "21.23", "92.20", "6.39", "1.3566", 3),
    @("515220", "国泰中证煤炭ETF", "32.88", "99.67", "3.51", "1.1541", 10),
    @("213001", "宝盈鸿利收益灵活配置混合A", "17.98", "90.37", "4.28", "0.7695", 5),
    @("001556", "天弘中证500指数增强A", "41.41", "94.29", "1.77", "0.7330", 6),
    @("013233", "华夏中证500指数智选增强A", "39.48", "92.73", "1.54", "0.6080", 4),
    @("007994", "华夏中证500指数增强A", "31.45", "92.72", "1.54", "0.4843", 4),
    @("006682", "景顺长城中证500指数增强", "16.63", "87.75", "1.90", "0.3160", 9),
    @("001557", "天弘中证500指数增强C", "13.97", "94.29", "1.77", "0.2473", 6),
    @("010751", "宝盈优质成长混合A", "5.64", "92.80", "4.28", "0.2414", 9),
    @("001543", "宝盈新锐灵活配置混合A", "3.21", "93.26", "4.88", "0.1566", 6),
    @("005062", "博时中证500指数增强A", "6.14", "90.01", "1.57", "0.0964", 6),
    @("007995", "华夏中证500指数增强C", "5.45", "92.72", "1.54", "0.0839", 4),
    @("006440", "中信建投中证500指数增强A", "5.78", "94.71", "1.22", "0.0705", 1),
    @("013234", "华夏中证500指数智选增强C", "4.28", "92.73", "1.54", "0.0659", 4),
    @("003016", "中金中证500指数增强A", "4.71", "93.78", "1.25", "0.0589", 6),
    @("162216", "泰达宏利中证500指数增强（LOF）", "4.44", "93.67", "1.23", "0.0546", 9),
    @("970041", "国海量化优选一年持有股票A", "7.70", "91.93", "0.67", "0.0516", 10),
    @("510170", "国联安上证大宗商品股票ETF", "2.22", "98.14", "2.24", "0.0497", 9),
    @("970042", "国海量化优选一年持有股票C", "6.95", "91.93", "0.67", "0.0466", 10),
    @("006441", "中信建投中证500指数增强C", "3.11", "94.71", "1.22", "0.0379", 1),
    @("010752", "宝盈优质成长混合C", "0.78", "92.80", "4.28", "0.0334", 9),
    @("007581", "宝盈鸿利收益灵活配置混合C", "0.73", "90.37", "4.28", "0.0312", 5),
    @("011824", "浙商汇金量化臻选股票型证券投资基金A", "1.54", "92.80", "1.51", "0.0233", 7),
    @("519097", "新华中小市值优选混合", "0.75", "62.70", "3.04", "0.0228", 9),
    @("710301", "富安达增强收益债券A", "0.61", "20.20", "2.95", "0.0180", 2),
    @("003578", "中金中证500指数增强C", "1.44", "93.78", "1.25", "0.0180", 6),
    @("005795", "博时中证500指数增强C", "1.14", "90.01", "1.57", "0.0179", 6),
    @("006729", "万家中证500指数增强A", "1.04", "93.64", "1.26", "0.0131", 7),
    @("007578", "宝盈新锐灵活配置混合C", "0.20", "93.26", "4.88", "0.0098", 6),
    @("710302", "富安达增强收益债券C", "0.26", "20.20", "2.95", "0.0077", 2),
    @("006730", "万家中证500指数增强C", "0.61", "93.64", "1.26", "0.0077", 7),
    @("011825", "浙商汇金量化臻选股票型证券投资基金C", "0.47", "92.80", "1.51", "0.0071", 7),
    @("003717", "中银量化精选灵活配置混合A", "0.49", "90.38", "1.18", "0.0058", 6),
    @("970073", "东证融汇成长优选混合A", "0.68", "82.02", "0.84", "0.0057", 9),
    @("005260", "银华稳健增利灵活配置混合A", "0.32", "91.49", "1.01", "0.0032", 5),
    @("970074", "东证融汇成长优选混合C", "0.27", "82.02", "0.84", "0.0023", 9),
    @("515510", "嘉实中证500成长估值ETF", "0.15", "98.79", "1.20", "0.0018", 10),
    @("519117", "浦银安盛基本面400指数", "0.24", "92.63", "0.74", "0.0018", 3),
    @("162907", "泰信中证锐联基本面400指数（LOF）", "0.23", "94.61", "0.77", "0.0018", 3),
    @("006157", "财通量化核心优选混合", "0.09", "92.85", "1.44", "0.0013", 6),
    @("005261", "银华稳健增利灵活配置混合C", "0.02", "91.49", "1.01", "0.0002", 5),
    @("010484", "中银量化精选灵活配置混合C", "0.01", "90.38", "1.18", "0.0001", 6)
)

# Template sheet ("2021-Q4") already has the exact header/format we need
# for the fund-detail table; copy its formatting across.
$template = $wb.Worksheets.Item(3)
$template.Range("B1:H1").Copy()
$detail.Range("B1:H1").PasteSpecial(-4122)

$rowCount = $fundRows.Count
$lastRow = $rowCount + 1
$template.Range("A2:H2").Copy()
$detail.Range("A2:H$lastRow").PasteSpecial(-4122)

$detail.Cells.Item(1, 2).Value = "基金代码"
$detail.Cells.Item(1, 3).Value = "基金名称"
$detail.Cells.Item(1, 4).Value = "基金规模"
$detail.Cells.Item(1, 5).Value = "股票总仓位"
$detail.Cells.Item(1, 6).Value = "仓位占比"
$detail.Cells.Item(1, 7).Value = "持有市值(亿元)"
$detail.Cells.Item(1, 8).Value = "仓位排名"

# Columns B-G hold text that looks numeric (fund codes, percentages,
# etc.); force text format before assigning so they keep their original
# textual representation (leading zeros, trailing zeros, ...).
$textRange = $detail.Range("B2:G$lastRow")
$textRange.NumberFormat = "@"

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $rec = $fundRows[$i]
    $detail.Cells.Item($r, 1).Value = $i
    $detail.Cells.Item($r, 2).Value = $rec[0]
    $detail.Cells.Item($r, 3).Value = $rec[1]
    $detail.Cells.Item($r, 4).Value = $rec[2]
    $detail.Cells.Item($r, 5).Value = $rec[3]
    $detail.Cells.Item($r, 6).Value = $rec[4]
    $detail.Cells.Item($r, 7).Value = $rec[5]
    $detail.Cells.Item($r, 8).Value = $rec[6]
}

$textRange.ClearFormats()

# ---------------------------------------------------------------------
# Step 3: update the new "总计" roll-up sheet: add a 2022-Q1 row on top
# and keep the historical rows below it (shifted down by one).
# ---------------------------------------------------------------------
$total.Range("A2:D2").Copy()
$total.Range("A5:D5").PasteSpecial(-4122)

$summaryRows = @(
    @(0, "2022-Q1", 42, 6.92),
    @(1, "2021-Q4", 14, 4.98),
    @(2, "2021-Q3", 27, 8.81),
    @(3, "2021-Q2", 14, 1.79)
)
for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $rec = $summaryRows[$i]
    $total.Cells.Item($r, 1).Value = $rec[0]
    $total.Cells.Item($r, 2).Value = $rec[1]
    $total.Cells.Item($r, 3).Value = $rec[2]
    $total.Cells.Item($r, 4).Value = $rec[3]
}

# ---------------------------------------------------------------------
# Step 4: restore the original active sheet/tab selection (sheet 1).
# ---------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
